$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray <w:bookmarkStart/.../bookmarkEnd> pair (_GoBack)
#    that sits right after "Work experience:" at the end of that
#    paragraph's runs.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Locate the trailing "." run that immediately follows the
#    "https://mineplex.market/" hyperlink (bold + single-underline run)
#    and replace it with:
#      - ", " text
#      - a HYPERLINK field (begin / instrText / empty / separate)
#        pointing at https://mineplex.io/
#      - the field's displayed text "https://mineplex.io/" split
#        across several runs (matching the original mixed bold styling)
#      - the _GoBack bookmark re-inserted in the middle of "mineplex"
#      - the field's closing fldChar
# ---------------------------------------------------------------------
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$rng = $d.Content
$found = $rng.Find.Execute("https://mineplex.market/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not locate 'https://mineplex.market/' in the document"
}
# The hyperlink run is immediately followed (after a zero-width field/run
# boundary) by a single-character run holding the trailing ".". Scan
# forward from the end of the hyperlink match for the first position that
# actually yields visible text, and use that one character as our target.
$pos = $rng.End
while ($d.Range($pos, $pos + 1).Text.Length -eq 0) {
    $pos = $pos + 1
}
$periodRng = $d.Range($pos, $pos + 1)
if ($periodRng.Text -ne ".") {
    throw "expected trailing '.' run, found [$($periodRng.Text)] instead"
}

$xml = @"
<w:p xmlns:w='$w'>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:u w:val="single"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:t xml:space="preserve">, </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:u w:val="single"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:fldChar w:fldCharType="begin"/>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:u w:val="single"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:instrText xml:space="preserve"> HYPERLINK "https://mineplex.io/" </w:instrText>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:u w:val="single"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:u w:val="single"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:fldChar w:fldCharType="separate"/>
</w:r>
<w:r>
<w:rPr>
<w:rStyle w:val="a5"/>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:t>https</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rStyle w:val="a5"/>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:t>://</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rStyle w:val="a5"/>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:t>mine</w:t>
</w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r>
<w:rPr>
<w:rStyle w:val="a5"/>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:t>plex</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rStyle w:val="a5"/>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:t>.io/</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:sz w:val="28"/>
<w:szCs w:val="28"/>
<w:u w:val="single"/>
<w:lang w:val="en-US" w:eastAsia="ru-RU"/>
</w:rPr>
<w:fldChar w:fldCharType="end"/>
</w:r>
</w:p>
"@

$periodRng.InsertXML($xml)
